# Add additional training/test-set instances to the lintroller list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Iterative-calculation "maximum change" setting (calcPr iterateDelta: 0.001 -> 0.0001)
$excel.MaxChange = 0.0001

# New rows 22-31: lintroller_131 .. lintroller_140
$lowerNames = @(
    "lintroller_131",
    "lintroller_132",
    "lintroller_133",
    "lintroller_134",
    "lintroller_135",
    "lintroller_136",
    "lintroller_137",
    "lintroller_138",
    "lintroller_139",
    "lintroller_140"
)
for ($i = 0; $i -lt $lowerNames.Length; $i++) {
    $r = 22 + $i
    $ws.Cells.Item($r, 1).Value = $lowerNames[$i]
    $ws.Cells.Item($r, 2).Value = 1
}

# New rows 32-36: Lintroller_141 .. Lintroller_145
$upperNames = @(
    "Lintroller_141",
    "Lintroller_142",
    "Lintroller_143",
    "Lintroller_144",
    "Lintroller_145"
)
for ($i = 0; $i -lt $upperNames.Length; $i++) {
    $r = 32 + $i
    $ws.Cells.Item($r, 1).Value = $upperNames[$i]
    $ws.Cells.Item($r, 2).Value = 1
}

# Move the active selection to E11, matching the saved view state.
$ws.Range("E11").Select() | Out-Null
